$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal de travail")

# Row 58: hours 4 -> 5, comment gets extra text appended
$ws.Range("C58").Value = 5
$ws.Range("D58").Value = "Backend: Validation des request body, reset config default, fixs"

# Row 59: new journal entry (was an empty styled row)
$ws.Range("A59").Value = 45101
$ws.Range("B59").Value = "Implémentation"
$ws.Range("C59").Value = 3
$ws.Range("D59").Value = "Frontend: refactor"

# Move the active selection from D60 to D61
[void]$ws.Range("D61").Select()
